# TNRSequencer.xlsx - insert "RO.ACT" into column A, row 3 of sheet "LIST",
# shifting the previous A3 value ("AD.SEC.002.FON.01") down to A4.
# Only column A is affected - other columns (D:H) in rows 3/4 stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LIST")

# Capture the value currently sitting in A3 before we overwrite it
$previousA3 = $ws.Range("A3").Value2

# Insert the new label at A3
$ws.Range("A3").Value = "RO.ACT"

# Move the old A3 content down into A4 (A4 was previously empty)
$ws.Range("A4").Value = $previousA3

# Update the selection to a single cell (A9) instead of the prior range (A8:A9)
$ws.Range("A9").Select()
